$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet from "My Series" to "Data"
$ws.Name = "Data"

# 2. Update the CEIC refresh-link comment payload on A1
$commentText = "FR0AAB+LCAAAAAAAAAOlWc1zI0cV/1emdIKDNDOSveuo2pPSh21USJbLkvE6F2o007Yaz0wr0z22ddsLFdiQE0mKkEDYUyiKDTlAFZWF8L+kkHc58S/w+mO+JHnXMlsu7/T76tevX7/36zZ69yYMjCscM0Kj3YpdsyoGjjzqk+hit5Lw86r9qPKug/ZuPBwcubEbYg7CBmhFrHnDyG5lyvmsaZrX19e160aNxhdm3bJs88mgP/KmOHSrJGLcjTxcybT8t2tVHNTxwwHmru9yV2nuVnqjXq2DidcF2sCN3Asc19oJIxFmbC/ihBPMhGaMXY473cFP1MKceu1RzUbmCj2XbCck8JVcSVLRtRxMi8ckxE7dsneq1k61YY9tu7nVaFp2rfFo571UMRNEfZfxEY6viCcJI+6GM6lu7TRs295qWDDbWiGwlQfAQcPAP8ZXhGG/g4OAbRQRU29gy+Ow6s2CaSGzoKsNPdyFg9idTceEB3gzN9rdsRFG2pfciIP2aYw9iN+DXDrE18NYh3U86wN3PCUxn3fd+ca2ThiOhzMRpM1UHdSlEW8FOOYnM9hr7EMqAMPhcYKReQczV+oS5sE3iRLsO+duwIpKJSY6pfElm7kePoRzbAob11FAXR8SjhPGiZdPusJARzGdgUWYvE0Dfx+sauE1jMxyL4IQi2nblF7m3q1jIrmrcn9hT0OXp+IrdDSa0uthFMxHyYR5MZlgv9tOpdfykDiQWruTME5D8CInIUUrUAYDcw7/4Awuc1AXeyR0g6MA4sicLTBUIqBWwuk54R0aJGHEUreWqOgUFjXGN9kiszEawv5GIu406kWpvIr0WlZZ4ZheZ3OuMmQcCuQW89IdX2UsC3eBlu7gKkduiljlPgmgRRS3o0AtJ8ZoijFfmxWKg0Q13BdNx2nPD5NwAidsAsfsSs7KkJnzEaQqpDv45VjQSKryZ2xZTfkDfmRstBf5d8ulTATTFeZybOAtkRCsKWgHbnQJ1FPCp4etdC1rOEhF4E75VR6CwzsL3LkkZ1Eq0lAv8oLEx6om9KJzmaLCN7Wpd7LRCqkPx9xBbjQfz2dQmhlpcvjYrUCzbjIeAxyoOB5NIh7PRfFAphZ9mw5LJpGcwA3urXMe4/cTQCHz/STyOtS//2y+is5JRPj9PaRJrCri/VVk9ERxTFgXizIj6/699b1N1sTijcTDCIc0It79ow1BFt77D1gIS0/VvTWwOl/3lg+gtavOJ876vdVigJDQ6zaapsUY9YhMVn08/IK+eceR6eJzNwkAvnHoshdZ7V0moxa7XJYpktBJHKQV0BHgmAE69vyw5gF+EAiw5tFQEEwApacjZBblBQjy8F500XejiwRgRlZXlulZ/RUtchy7ERPLyVDFUileL4TSOqXQjqOK1zCRiaCKFwUuMpfk0BiHMxq7wQACQ/Z12mnIBGhk4PKpHkFvC7CXBtnMVTOtsmep428Tk01KLUMceF0ml4hSSKxF4fBcJqchscoBHMug4wZkEquqmrbydTzYsBwfpvVXLG5DrJjuAdzFoPv+GM8FQM8Hmi5T1k4ZKoFFIXVGx1s79W3r0Q4gGzFGcsUHFG6AUQiiRhdPeNPo0hAL2Nc02u2mUWLPKCMcNleFSrejB1sp6qP9tO5L1waQviVCmQ/w44JAV1mRyxi5vDMAUDoN5gVJtfI+9UDw9uN/Lb79YvH028Wz5//9x28Xn/9z8fOv4KPdhl+a+eI3t3/+Tq1aKaGxOwmwdGzc3tmxGoAEcxISMTclaPYTj0va2ZnEytkY6SufHHT2ep2DfluWmYyYqqtOY4rb5Jwm+XCkFiMnkjttpgmiRJxxWrb0uMQtdC5H3O6ucFm6yL9LUcXi1cuvXr18cae2DlgOweCqX69aW29FaHBPblTt+h0IrV9qDUJ4W9zJ7UZBeEkGHau+kMWp5zsN23rHevxY3MMzWprX64SWWdrS2L0wl/QUqaNQU5YCxXHKlPk/hsOSsdWJKAx0qv7tw9d/+aQkpaOrKWUr4JxEN2IyMx1I04fHY2M0PDnu7BnjvZHIk5xXkFPG3yCsZ8/OVTFpVk9dKaXWUZW5Nyouz3gQ02S2UjFy6hrJtbVjlbukKWOz6lrOWyOvd+2LX6xT0Cvp5qA1e+Yo0lCJo0gFvj6Cn333779/sHjx5eLZH24//OXi+ceLj56X7OjZsjcASF04IMVhlslQxXRnWaKg05GM6aX100In0URxZTqiJOLMse1teV3SQwS6tjAn/0e9ELqbtCzDBvQlCvqRy/ZuuD6sziEyywRwdOZCY6X5NTMjqLqch/c/v/v97ed/vf30m9cf/Gnx7I+LX3366uWXr79+rk7S7Sff3H70ta7cy8Vd+iIurwrvGfI1xDPECTNEmza+f/prI6LcAHRhJLLKfP/0s4Ix4ajEIbllQG+ZI2UXVkSLykLPKLiS+VDSy1RUr++IttTIJHRjojPi5ZO8VxWmAIIakvGD3riaMGxQAE4/hJWUhXPl++ppFdUmjx5bdbuuucobsYSJywqhPwjoxA2MlCHfGpZESlpvVshl5XwH/WG71c9FlBPD2MexSEP1gVL0KNpEj6WjNNUKFOACxvOSQDwPrYitsjLLhXJm6neW85YvHt3Xv0qUJFAniWOFdiL9cj9KZoB70we5u/nykbIAdQ8VLC2C33zc65b5MC5wobmV2YIg+bJCaZaqVj0mnnQUcj0UocmHwCs9bEI49OO8Qk9XBLbEFIVnL45pvLb65JxUbACgGSqKmUc8k5F7qgC2n+9VSkgrXvahLnl6hbSLA8w3e7k2c+0BIOCH6sLeb6raY8PA18Hc7JaRhSU3UHy+F4ny/77eq2RrxTGAJfHWt/Fze3pHPYar7YbeqKVIRXHZg9n1e/g+iRl/IiqB/lKUs4xyplDnE0f/0eaJGp85jW1FAAGzaN0suZkeXa7+iEKDPgnJhjdAKz3fZSMQy9lMwbLeZpkiWsshvgHQWLAARXHyM2gb6slkE2sqYaGWZvrimZKRiynf1LHHExf7eGJVvQmuV7d8APfvYAw3Aht+u14drhHb4pFTG4fKQfD1hpOY6Yblf9x0/gfhIDnFFR0AAA=="
[void]$ws.Range("A1").Comment.Text($commentText)

# 3. Relabel the "Function Description" row to "Function Information"
$ws.Range("A11").Value = "Function Information"

# 4. Update the "Last Update Time" serial date value
$ws.Range("B14").Value = 42229

# 5. Strip the bold/red highlight font from the B27:B36 data column and
#    switch its custom number format from "0.0000" to "###0.0000"
$dataRange = $ws.Range("B27:B36")
$dataRange.ClearFormats()
$dataRange.NumberFormat = "###0.0000"

